# Apply commit "otvet na 1 vopros tatjani":
#   - append a new run "wahmatnij" to the end of the first question's
#     paragraph (same rStyle/rFonts/color as the trailing space run that
#     already ends that paragraph);
#   - the document's "_GoBack" bookmark (which Word automatically drops at
#     the point of the most recent edit) now sits right after that new
#     run instead of at the end of the document.

$d = $word.ActiveDocument

$newWord = "wahmatnij"

# --- Step 1: append the new run at the end of paragraph 1 -----------------
$p1 = $d.Paragraphs.Item(1)
$insertPoint = $d.Range($p1.Range.End - 1, $p1.Range.End - 1)
$insertPoint.InsertAfter($newWord)

# Re-fetch paragraph 1 (positions shifted after the insert) and grab the
# range that now covers the text we just inserted so we can format it to
# match the neighbouring "apple-converted-space" run.
$p1 = $d.Paragraphs.Item(1)
$newTextEnd = $p1.Range.End - 1
$newRunStart = $newTextEnd - $newWord.Length
$newRun = $d.Range($newRunStart, $newTextEnd)
$newRun.Style = "apple-converted-space"
$newRun.Font.NameAscii = "Arial"
$newRun.Font.Name = "Arial"
$newRun.Font.NameBi = "Arial"
$newRun.Font.Color = 0

# --- Step 2: move the "_GoBack" bookmark to the end of paragraph 1 --------
# Bookmarks.Add collapses to position 0 when handed a zero-length range in
# this host, so bracket a throw-away marker character, bookmark around it,
# then delete the marker - leaving the bookmark collapsed exactly where we
# want it. Adding a bookmark named "_GoBack" also removes the old one
# (bookmark names are unique), so this both relocates and keeps a single
# pair in the document.
$p1 = $d.Paragraphs.Item(1)
$markerInsert = $d.Range($p1.Range.End - 1, $p1.Range.End - 1)
$markerInsert.InsertAfter("X")

$p1 = $d.Paragraphs.Item(1)
$markerRange = $d.Range($p1.Range.End - 2, $p1.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $markerRange)

$p1 = $d.Paragraphs.Item(1)
$markerRange = $d.Range($p1.Range.End - 2, $p1.Range.End - 1)
$markerRange.Delete()

Write-Output "Inserted run and relocated _GoBack bookmark."
